$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: force a clean run-boundary split at a collapsed Range point.
# Adding then immediately deleting a bookmark at a point forces the
# serializer to keep that point as a run boundary, without leaving any
# bookmark residue behind in the saved OOXML.
# ------------------------------------------------------------------
function Split-RunAt($point) {
    $rng = $d.Range($point, $point)
    $bm = $d.Bookmarks.Add("zzsplit", $rng)
    $d.Bookmarks("zzsplit").Delete()
}

# ------------------------------------------------------------------
# Step 1: apply BOTH content edits first. (Any text-changing edit to a
# paragraph re-flows/merges its runs, so all content changes must land
# before we carve out the final run boundaries below.)
# ------------------------------------------------------------------

# 1a. Insert the new sentence right after "...enter the gauntlet. " and
#     before "He also estimates...".
$r1 = $d.Content
$r1.Find.Execute("enter the gauntlet. He also", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insStart = $r1.Start + [string]"enter the gauntlet. ".Length
$insPt = $d.Range($insStart, $insStart)
$insPt.InsertBefore("I just shifted the whole distribution earlier by a week to cover this. ")

# 1b. Swap the old single trailing space (after "...distinct on/off
#     dates?") for the new " Or just a poorly fitting normal. " text.
$r2 = $d.Content
$r2.Find.Execute("dates? Mostly seal predation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spaceStart = $r2.Start + [string]"dates?".Length
$rOldSpace = $d.Range($spaceStart, $spaceStart + 1)
$rOldSpace.Text = " Or just a poorly fitting normal. "

# ------------------------------------------------------------------
# Step 2: now that the paragraph's final text is in place, carve out
# the run boundaries that match the target structure.
# ------------------------------------------------------------------

# Boundary before "I just shifted..." / after "...He also..." sentence.
$rA = $d.Content
$rA.Find.Execute("enter the gauntlet. I just shifted", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryA = $rA.Start + [string]"enter the gauntlet. ".Length
Split-RunAt $boundaryA

$rB = $d.Content
$rB.Find.Execute("cover this. He also estimates", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryB = $rB.Start + [string]"cover this. ".Length
Split-RunAt $boundaryB

# Boundary before the lone space run / after "...good starting point."
# and the boundary right after that space, before "This is not...".
$rE = $d.Content
$rE.Find.Execute("good starting point. This is not", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryE = $rE.Start + [string]"good starting point.".Length
Split-RunAt $boundaryE
Split-RunAt ($boundaryE + 1)

# Boundaries around the new " Or just a poorly fitting normal. " run.
$rC = $d.Content
$rC.Find.Execute("distinct on/off dates? Or just", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryC = $rC.Start + [string]"distinct on/off dates?".Length
Split-RunAt $boundaryC

$rD = $d.Content
$rD.Find.Execute("poorly fitting normal. Mostly seal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryD = $rD.Start + [string]"poorly fitting normal.".Length
Split-RunAt $boundaryD
Split-RunAt ($boundaryD + 1)
